$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(45826,45828,45831,45832,45833,45834,45835,45838,45839,45840,45841,45845,45846,45847,45848,45849,45852,45853,45854,45855,45856,45859,45860,45861,45862,45863,45866,45867,45868,45869,45870,45873,45874,45875,45876,45877,45880,45881,45882,45883,45884,45887,45888,45889,45890,45891,45894,45895,45896,45897,45898,45902,45903,45904,45905,45908,45909,45910,45911,45912,45915)

$startRow = 1321
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = 0
}
